$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new subscriber (+5516999928423, DDD 16, 2024-10-15) was added at the top
# of the list, pushing every existing data row down by one.
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above it (the header);
# re-copy the formatting used by the rest of the data rows instead.
$ws.Range("A3:C3").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)  # xlPasteFormats

# Force text so the phone number keeps its "+" prefix, "16" isn't coerced to
# a number, and the date stays a literal string like the rest of the column.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 3).NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "+5516999928423"
$ws.Cells.Item(2, 2).Value = "16"
$ws.Cells.Item(2, 3).Value = "2024-10-15"

# Restore General formatting to match the other data rows now that the
# values are stored as text.
$ws.Cells.Item(2, 1).NumberFormat = "General"
$ws.Cells.Item(2, 2).NumberFormat = "General"
$ws.Cells.Item(2, 3).NumberFormat = "General"
